$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet is protected; temporarily unprotect to edit the locked cells,
# then re-apply protection afterwards.
$ws.Unprotect()

# Update the confidentiality/date notice text on row 11 (A11)
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-13 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for each holding row
$ws.Range("D2").Value = 0.5306642001721305
$ws.Range("E2").Value = -0.006907652068942394

$ws.Range("D3").Value = 0.2677408637935658
$ws.Range("E3").Value = -0.0002672367717799773

$ws.Range("D4").Value = 0.05005160372143046
$ws.Range("E4").Value = -0.01228501228501233

$ws.Range("D5").Value = 0.09478613289965959
$ws.Range("E5").Value = -0.01773919607473118

$ws.Range("D6").Value = 0.02721814483363739
$ws.Range("E6").Value = -0.02050987157370143

$ws.Range("D7").Value = 0.02953905457957641
$ws.Range("E7").Value = -0.01710695632427539

$ws.Range("E8").Value = -0.00709707219918676

# Re-apply sheet protection (same settings as before the edit).
$ws.Protect()
